# Get rid of "All" value in Tech Spec: rename "total_staff" sheet to "total_staff_by_type"
# and update the active selection on that sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("total_staff")
$ws.Name = "total_staff_by_type"

$ws.Activate()
$ws.Range("D17").Select()
